$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9906358122825623
$ws.Range("B1").Value = 2.211426496505737
$ws.Range("C1").Value = 4.972152709960938
$ws.Range("D1").Value = 1.736373782157898
$ws.Range("E1").Value = 1.295992970466614
